$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two old test case rows (rows 3 and 4)
$ws.Rows("3:4").Delete()

# Replace the remaining test case (row 2) with the new iAuthor test case
$ws.Range("A2").Value = "iAU_TC_ID_107"
$ws.Range("B2").Value = "@RegressionA Validation of Blueprints list page"
$ws.Range("C2").Value = "passed"
